$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 619 (shifts existing rows 619:686 down to 620:687)
$ws.Rows.Item(619).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A619").Value = 3
$ws.Range("B619").Value = "Femacal de La Calera"
$ws.Range("C619").Value = "Coquimbo"
$ws.Range("D619").Value = 45212
$ws.Range("E619").Value = 5
$ws.Range("F619").Value = 100112031
$ws.Range("G619").Value = "Poroto verde"
$ws.Range("H619").Value = "Sin especificar"
$ws.Range("I619").Value = "Primera"
$ws.Range("J619").Value = 38
$ws.Range("K619").Value = 34000
$ws.Range("L619").Value = 34000
$ws.Range("M619").Value = 34000
$ws.Range("N619").Value = "`$/malla 25 kilos"
$ws.Range("O619").Value = "Provincia de Limarí"
$ws.Range("P619").Value = 1360
$ws.Range("Q619").Value = 25
$ws.Range("R619").Value = "Hortaliza"
